$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F, matching style of existing header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Timestamp values for the data rows F2:F12
$timestamps = @(
    "2021-10-05 10:52:09.153651",
    "2021-10-05 10:52:09.153663",
    "2021-10-05 10:52:09.153666",
    "2021-10-05 10:52:09.153670",
    "2021-10-05 10:52:09.153673",
    "2021-10-05 10:52:09.153677",
    "2021-10-05 10:52:09.153680",
    "2021-10-05 10:52:09.153683",
    "2021-10-05 10:52:09.153686",
    "2021-10-05 10:52:09.153689",
    "2021-10-05 10:52:09.153692"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
